$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits (title rows for the two tables) ---
# Row 14: "Share of Families With/with Student Loan Debt..." - hyphens -> en dashes, "With" -> "with"
$ws.Range("A14").Value2 = "Share of Families with Student Loan Debt for Those Ages 25–55, 1989–2016"

# Row 1: "Average Family Student Loan Debt..." - hyphens -> en dashes
$ws.Range("A1").Value2 = "Average Family Student Loan Debt for Those Age 25–55, 1989–2016"

# --- Formatting edits ---
# Row 1 title: align to top and drop the old custom (taller) row height
$ws.Range("A1").VerticalAlignment = -4160
$ws.Rows.Item(1).AutoFit()

# Percent tables (rows 16-25, cols B:D): right-align and top-align instead of the old
# center alignment
$ws.Range("B16:D25").HorizontalAlignment = -4152
$ws.Range("B16:D25").VerticalAlignment = -4160

# Move the active selection to F14 and drop the old scrolled-down view (A21 top-left)
$ws.Range("F14").Select()
